$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1047.7273
$ws.Range("I32").Value = 1899
$ws.Range("J32").Value = 962.6
$ws.Range("K32").Value = 1899
$ws.Range("L32").Value = 962.6
$ws.Range("M32").Value = -1573
$ws.Range("N32").Value = -1614.6
$ws.Range("H38").Value = 183178.8
$ws.Range("I38").Value = 183178.8
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 549536.3999999999
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -549164.3999999999
$ws.Range("H40").Value = 7866.2666
$ws.Range("I40").Value = 6000
$ws.Range("K40").Value = 6000
$ws.Range("M40").Value = -5825
$ws.Range("H52").Value = 28856.172
$ws.Range("I52").Value = 500049.5
$ws.Range("J52").Value = 299
$ws.Range("K52").Value = 1500148.5
$ws.Range("L52").Value = 897
$ws.Range("M52").Value = -1499988.5
$ws.Range("N52").Value = -1217
$ws.Range("H53").Value = 27778568
$ws.Range("I53").Value = 66667820
$ws.Range("J53").Value = 531.1429000000001
$ws.Range("K53").Value = 66667820
$ws.Range("L53").Value = 531.1429000000001
$ws.Range("M53").Value = -66667183
$ws.Range("N53").Value = -1805.1429
$ws.Range("H104").Value = 591.2857
$ws.Range("I104").Value = 591.2857
$ws.Range("K104").Value = 1773.8571
$ws.Range("M104").Value = -26.85710000000017
$ws.Range("H132").Value = 1752.8372
$ws.Range("I132").Value = 1637.5897
$ws.Range("J132").Value = 2876.5
$ws.Range("K132").Value = 4912.7691
$ws.Range("L132").Value = 8629.5
$ws.Range("M132").Value = -2382.7691
$ws.Range("N132").Value = -13689.5
$ws.Range("H137").Value = 4677.7036
$ws.Range("I137").Value = 4679.0835
$ws.Range("K137").Value = 14037.2505
$ws.Range("M137").Value = -11487.2505
$ws.Range("H138").Value = 6905.5
$ws.Range("J138").Value = 7961.273
$ws.Range("L138").Value = 23883.819
$ws.Range("N138").Value = -34163.819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3033.8057
$ws.Range("I132").Value = 3164.182
$ws.Range("K132").Value = 9492.545999999998
$ws.Range("M132").Value = -6962.545999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4418.2
$ws.Range("I20").Value = 4397.6
$ws.Range("K20").Value = 4397.6
$ws.Range("M20").Value = -4150.6
$ws.Range("H57").Value = 79998.664
$ws.Range("J57").Value = 79998.664
$ws.Range("L57").Value = 79998.664
$ws.Range("N57").Value = -81438.664
$ws.Range("H58").Value = 117999.664
$ws.Range("J58").Value = 117999.664
$ws.Range("L58").Value = 117999.664
$ws.Range("N58").Value = -118587.664
$ws.Range("H59").Value = 91931.5
$ws.Range("J59").Value = 91931.5
$ws.Range("L59").Value = 91931.5
$ws.Range("N59").Value = -93625.5
$ws.Range("H60").Value = 48789.832
$ws.Range("J60").Value = 48789.832
$ws.Range("L60").Value = 48789.832
$ws.Range("N60").Value = -49987.832
$ws.Range("H86").Value = 2836024.8
$ws.Range("I86").Value = 3402229.8
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 3402229.8
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -3401106.8
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 2836024.8
$ws.Range("I89").Value = 3402229.8
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 17011149
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -17005533
$ws.Range("N89").Value = -36232
$ws.Range("H124").Value = 89499.25
$ws.Range("J124").Value = 89499.25
$ws.Range("L124").Value = 89499.25
$ws.Range("N124").Value = -99319.25
$ws.Range("H132").Value = 49750
$ws.Range("J132").Value = 49750
$ws.Range("L132").Value = 49750
$ws.Range("N132").Value = -59870
$ws.Range("H133").Value = 47388.832
$ws.Range("H134").Value = 47514
$ws.Range("I134").Value = 3988.8235
$ws.Range("J134").Value = 170835.33
$ws.Range("K134").Value = 11966.4705
$ws.Range("L134").Value = 512505.99
$ws.Range("M134").Value = -9431.470499999999
$ws.Range("N134").Value = -517575.99
$ws.Range("H136").Value = 79998.664
$ws.Range("J136").Value = 79998.664
$ws.Range("L136").Value = 79998.664
$ws.Range("N136").Value = -90198.664
$ws.Range("H137").Value = 75000
$ws.Range("J137").Value = 75000
$ws.Range("L137").Value = 75000
$ws.Range("N137").Value = -85200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1196.7576
$ws.Range("I132").Value = 1115.7
$ws.Range("K132").Value = 3347.1
$ws.Range("M132").Value = -817.1000000000004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 789215.3
$ws.Range("I132").Value = 126057.5
$ws.Range("K132").Value = 1134517.5
$ws.Range("M132").Value = -1131987.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5234.7
$ws.Range("I122").Value = 2899.6667
$ws.Range("J122").Value = 6235.4287
$ws.Range("K122").Value = 8699.000100000001
$ws.Range("L122").Value = 18706.2861
$ws.Range("M122").Value = -6249.000100000001
$ws.Range("N122").Value = -23606.2861
$ws.Range("H132").Value = 86287.84
$ws.Range("I132").Value = 8674.700000000001
$ws.Range("K132").Value = 26024.1
$ws.Range("M132").Value = -23494.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6715.3687
$ws.Range("I7").Value = 7241
$ws.Range("K7").Value = 7241
$ws.Range("M7").Value = -7129
$ws.Range("H16").Value = 568.8570999999999
$ws.Range("I16").Value = 568.8570999999999
$ws.Range("K16").Value = 568.8570999999999
$ws.Range("M16").Value = -398.8570999999999
$ws.Range("H126").Value = 6715.3687
$ws.Range("I126").Value = 7241
$ws.Range("K126").Value = 21723
$ws.Range("M126").Value = -19253
$ws.Range("H132").Value = 7003.05
$ws.Range("I132").Value = 5716.5713
$ws.Range("K132").Value = 17149.7139
$ws.Range("M132").Value = -14619.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 660.4666999999999
$ws.Range("I107").Value = 703.0833
$ws.Range("J107").Value = 490
$ws.Range("K107").Value = 2109.2499
$ws.Range("L107").Value = 1470
$ws.Range("M107").Value = -189.2498999999998
$ws.Range("N107").Value = -5310
$ws.Range("H132").Value = 18582.492
$ws.Range("I132").Value = 1709.5745
$ws.Range("J132").Value = 75227.28999999999
$ws.Range("K132").Value = 5128.7235
$ws.Range("L132").Value = 225681.87
$ws.Range("M132").Value = -2598.7235
$ws.Range("N132").Value = -230741.87
$ws.Range("H136").Value = 9100395
$ws.Range("I136").Value = 11090307
$ws.Range("K136").Value = 33270921
$ws.Range("M136").Value = -33268371
